# Update cryptos list Price (D) and Volume(1h) (E) columns with refreshed scrape data
# Leading apostrophe forces text storage (matches original string cell type);
# Style is reset to Normal afterward so no stray number-format style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.939.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -5.39%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.821.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -4.43%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.66%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''329.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -2.76%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = '''  -0.49%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.4626'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -2.83%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.3848'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -3.83%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''45.87'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -3.00%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.07841'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.37%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.9595'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -3.22%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''21.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -6.35%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.796.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -5.16%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''5.632'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -4.88%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''6.847'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -3.73%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.06854'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.30%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -0.67%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''86.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.79%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.000009931'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -2.76%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''16.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -4.05%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.56%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''27.969.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = '''5.312'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -3.51%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''10.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -5.55%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''2.103'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -2.17%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''2.062.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -2.89%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''152.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -3.12%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''19.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.08%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''5.701'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -12.21%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -4.23%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''116.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -2.16%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''0.9360'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -5.92%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.09247'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -2.88%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''5.266'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -3.81%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''3.417'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -3.54%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -5.32%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -8.21%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.02146'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -4.25%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.145'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -4.05%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -0.60%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''7.576'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -2.22%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.5572'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -4.21%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''9.894'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -5.89%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.1766'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -3.01%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''1.242'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.22%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.212'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -9.96%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''11.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -4.76%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.5241'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -4.36%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.06993'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -5.71%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -6.47%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''112.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -2.92%  '
$ws.Range("E51").Style = "Normal"
